$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.006403796807923356
$ws.Range("D2").Value = 0.2301108526566651
$ws.Range("E2").Value = 0.1733927365998404
$ws.Range("F2").Value = 1.119256496137453
$ws.Range("G2").Value = 0.6172754033869978
$ws.Range("H2").Value = 0.6451706849114487
$ws.Range("I2").Value = 0.5624574055654961
$ws.Range("J2").Value = 0.1800749562410715
$ws.Range("O2").Value = 2.506829514151605

$ws.Range("C3").Value = 0.005607891203307958
$ws.Range("D3").Value = 0.2288708407612177
$ws.Range("E3").Value = 0.1704818582525682
$ws.Range("F3").Value = 1.086219363513223
$ws.Range("G3").Value = 0.5855465827597044
$ws.Range("H3").Value = 0.6344879103896801
$ws.Range("I3").Value = 0.5403325584106256
$ws.Range("J3").Value = 0.1749525151810474
$ws.Range("O3").Value = 2.416823902178777

$ws.Range("C4").Value = 0.005117407016410169
$ws.Range("D4").Value = 0.228207670226908
$ws.Range("E4").Value = 0.1687872790564455
$ws.Range("F4").Value = 1.066580177476183
$ws.Range("G4").Value = 0.5664367015796472
$ws.Range("H4").Value = 0.6282672455485283
$ws.Range("I4").Value = 0.5270823313488435
$ws.Range("J4").Value = 0.1719137970433167
$ws.Range("O4").Value = 2.363042203959708

$ws.Range("C5").Value = 0.004917077860174146
$ws.Range("D5").Value = 0.2279621468697002
$ws.Range("E5").Value = 0.1681200324150822
$ws.Range("F5").Value = 1.058739344460022
$ws.Range("G5").Value = 0.5587425543045441
$ws.Range("H5").Value = 0.6258174501459877
$ws.Range("I5").Value = 0.5217667888828004
$ws.Range("J5").Value = 0.1707022290642115
$ws.Range("O5").Value = 2.341498059993597

$ws.Range("C6").Value = 0.004883785960334563
$ws.Range("D6").Value = 0.2279228720738118
$ws.Range("E6").Value = 0.1680106442374765
$ws.Range("F6").Value = 1.0574471799641
$ws.Range("G6").Value = 0.557470575568999
$ws.Range("H6").Value = 0.6254158088869985
$ws.Range("I6").Value = 0.5208892200425694
$ws.Range("J6").Value = 0.1705026626183823
$ws.Range("O6").Value = 2.337943149268625

$ws.Range("C7").Value = 0.005114707143636821
$ws.Range("D7").Value = 0.2282042588655742
$ws.Range("E7").Value = 0.168778185967156
$ws.Range("F7").Value = 1.066473776248671
$ws.Range("G7").Value = 0.5663325581709557
$ws.Range("H7").Value = 0.6282338618372876
$ws.Range("I7").Value = 0.5270103039531975
$ws.Range("J7").Value = 0.1718973492214744
$ws.Range("O7").Value = 2.36275014500913

$ws.Range("C8").Value = 0.006129739967747128
$ws.Range("D8").Value = 0.229662932052193
$ws.Range("E8").Value = 0.1723698189763212
$ws.Range("F8").Value = 1.107731204184631
$ws.Range("G8").Value = 0.6062580231741208
$ws.Range("H8").Value = 0.6414170012355527
$ws.Range("I8").Value = 0.5547592695329868
$ws.Range("J8").Value = 0.1782866114388284
$ws.Range("O8").Value = 2.475487718664994

$ws.Range("C9").Value = 0.00810613374337521
$ws.Range("D9").Value = 0.2333016362271678
$ws.Range("E9").Value = 0.1801494517381954
$ws.Range("F9").Value = 1.19377010631348
$ws.Range("G9").Value = 0.6875148270817135
$ws.Range("H9").Value = 0.6699565692537988
$ws.Range("I9").Value = 0.6118365369392507
$ws.Range("J9").Value = 0.1916635588981563
$ws.Range("O9").Value = 2.708356419458369

$ws.Range("C10").Value = 0.00954995676701742
$ws.Range("D10").Value = 0.2364487660970838
$ws.Range("E10").Value = 0.1863161131496156
$ws.Range("F10").Value = 1.260133648751889
$ws.Range("G10").Value = 0.7490460878953797
$ws.Range("H10").Value = 0.6925673311477283
$ws.Range("I10").Value = 0.6554091411254745
$ws.Range("J10").Value = 0.2020134014545647
$ws.Range("O10").Value = 2.886696899377966

$ws.Range("C11").Value = 0.010205083468243
$ws.Range("D11").Value = 0.2379832533827795
$ws.Range("E11").Value = 0.1892198975436727
$ws.Range("F11").Value = 1.291013430937269
$ws.Range("G11").Value = 0.7774419290807657
$ws.Range("H11").Value = 0.7032114214600824
$ws.Range("I11").Value = 0.6755907234695968
$ws.Range("J11").Value = 0.2068362859216393
$ws.Range("O11").Value = 2.969418359926237

$ws.Range("C12").Value = 0.01045292532590736
$ws.Range("D12").Value = 0.2385790925782914
$ws.Range("E12").Value = 0.190333676672374
$ws.Range("F12").Value = 1.30280635855857
$ws.Range("G12").Value = 0.788253254152977
$ws.Range("H12").Value = 0.70729361931825
$ws.Range("I12").Value = 0.6832849390733458
$ws.Range("J12").Value = 0.2086791464770528
$ws.Range("O12").Value = 3.000972727489909

$ws.Range("C13").Value = 0.01039955883150157
$ws.Range("D13").Value = 0.23845011181335
$ws.Range("E13").Value = 0.1900931736462965
$ws.Range("F13").Value = 1.300262116966152
$ws.Range("G13").Value = 0.7859222396496364
$ws.Range("H13").Value = 0.7064121542274222
$ws.Range("I13").Value = 0.6816255428247189
$ws.Range("J13").Value = 0.208281517167265
$ws.Range("O13").Value = 2.994166714265816

$ws.Range("C14").Value = 0.01022547838216781
$ws.Range("D14").Value = 0.2380319777263509
$ws.Range("E14").Value = 0.1893112446014413
$ws.Range("F14").Value = 1.291981648967251
$ws.Range("G14").Value = 0.7783302109860415
$ws.Range("H14").Value = 0.7035462340550964
$ws.Range("I14").Value = 0.6762226907580526
$ws.Range("J14").Value = 0.206987567436812
$ws.Range("O14").Value = 2.97200975396953

$ws.Range("C15").Value = 0.0101188177875926
$ws.Range("D15").Value = 0.2377777802144863
$ws.Range("E15").Value = 0.1888341372831306
$ws.Range("F15").Value = 1.286922574371943
$ws.Range("G15").Value = 0.7736874903794444
$ws.Range("H15").Value = 0.7017974830692708
$ws.Range("I15").Value = 0.672920048056767
$ws.Range("J15").Value = 0.2061971411048233
$ws.Range("O15").Value = 2.958467881712124

$ws.Range("C16").Value = 0.009507109026550609
$ws.Range("D16").Value = 0.2363505509287052
$ws.Range("E16").Value = 0.186128327652284
$ws.Range("F16").Value = 1.258129489543307
$ws.Range("G16").Value = 0.7471985191324961
$ws.Range("H16").Value = 0.6918789240660033
$ws.Range("I16").Value = 0.6540974831537056
$ws.Range("J16").Value = 0.2017005264407459
$ws.Range("O16").Value = 2.881322956879842

$ws.Range("C17").Value = 0.009131416757433897
$ws.Range("D17").Value = 0.2355013134694701
$ws.Range("E17").Value = 0.1844936458195718
$ws.Range("F17").Value = 1.240642838012079
$ws.Range("G17").Value = 0.7310522582826877
$ws.Range("H17").Value = 0.6858859766670378
$ws.Range("I17").Value = 0.6426427683739604
$ws.Range("J17").Value = 0.1989714068215704
$ws.Range("O17").Value = 2.834405460557036

$ws.Range("C18").Value = 0.008915170754740132
$ws.Range("D18").Value = 0.2350225360989384
$ws.Range("E18").Value = 0.1835626944836406
$ws.Range("F18").Value = 1.230649981622719
$ws.Range("G18").Value = 0.7218034608710298
$ws.Range("H18").Value = 0.6824727215651478
$ws.Range("I18").Value = 0.6360882161767165
$ws.Range("J18").Value = 0.1974124816885876
$ws.Range("O18").Value = 2.807569719244043

$ws.Range("C19").Value = 0.008841926462061167
$ws.Range("D19").Value = 0.2348620939318664
$ws.Range("E19").Value = 0.1832490829698159
$ws.Range("F19").Value = 1.227277733275642
$ws.Range("G19").Value = 0.71867851532906
$ws.Range("H19").Value = 0.6813228456936429
$ws.Range("I19").Value = 0.6338747763727355
$ws.Range("J19").Value = 0.1968865085887757
$ws.Range("O19").Value = 2.798509342108332

$ws.Range("C20").Value = 0.009171426162701835
$ws.Range("D20").Value = 0.2355907144349771
$ws.Range("E20").Value = 0.1846667005584877
$ws.Range("F20").Value = 1.242497593576786
$ws.Range("G20").Value = 0.7327671101404576
$ws.Range("H20").Value = 0.6865204453715421
$ws.Range("I20").Value = 0.6438586334459728
$ws.Range("J20").Value = 0.1992608089660024
$ws.Range("O20").Value = 2.839384384731602

$ws.Range("C21").Value = 0.0102766165502004
$ws.Range("D21").Value = 0.2381543934124863
$ws.Range("E21").Value = 0.1895405312281824
$ws.Range("F21").Value = 1.294411123941288
$ws.Range("G21").Value = 0.7805585862729743
$ws.Range("H21").Value = 0.7043866263950349
$ws.Range("I21").Value = 0.6778082303904114
$ws.Range("J21").Value = 0.2073671826046279
$ws.Range("O21").Value = 2.978511557254819

$ws.Range("C22").Value = 0.01099752129009346
$ws.Range("D22").Value = 0.2399159389032235
$ws.Range("E22").Value = 0.1928085117554161
$ws.Range("F22").Value = 1.328919187989868
$ws.Range("G22").Value = 0.8121337845769574
$ws.Range("H22").Value = 0.7163634517996229
$ws.Range("I22").Value = 0.7002987435781591
$ws.Range("J22").Value = 0.2127615894428487
$ws.Range("O22").Value = 3.070777365247636

$ws.Range("C23").Value = 0.01061288889845002
$ws.Range("D23").Value = 0.2389679055618643
$ws.Range("E23").Value = 0.1910567634459568
$ws.Range("F23").Value = 1.310448518602939
$ws.Range("G23").Value = 0.795250265159126
$ws.Range("H23").Value = 0.709943727824367
$ws.Range("I23").Value = 0.6882674244892684
$ws.Range("J23").Value = 0.2098736547896891
$ws.Range("O23").Value = 3.02141079105553

$ws.Range("C24").Value = 0.009153338706475722
$ws.Range("D24").Value = 0.2355502667882092
$ws.Range("E24").Value = 0.1845884349490419
$ws.Range("F24").Value = 1.241658870283558
$ws.Range("G24").Value = 0.7319917200247232
$ws.Range("H24").Value = 0.6862335018541046
$ws.Range("I24").Value = 0.6433088445505177
$ws.Range("J24").Value = 0.1991299388575101
$ws.Range("O24").Value = 2.837132984039215

$ws.Range("C25").Value = 0.007572937895460541
$ws.Range("D25").Value = 0.2322340231463187
$ws.Range("E25").Value = 0.1779658061177258
$ws.Range("F25").Value = 1.169942548261687
$ws.Range("G25").Value = 0.6652128538278816
$ws.Range("H25").Value = 0.6619476978387411
$ws.Range("I25").Value = 0.5961090942242038
$ws.Range("J25").Value = 0.1879534732698005
$ws.Range("O25").Value = 2.644090394289606

